# The deck's theme (ppt/theme/theme1.xml, the theme driving the slide
# master / presentation) currently carries the "Integral" color scheme.
# The edit swaps the applied design's color scheme over to the stock
# "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), i.e. the
# classic default Office theme colors, while the font scheme / format
# scheme (already identical between the two themes in this deck) stay
# untouched.
#
# PowerPoint doesn't expose a raw RGB() literal helper in this host, so
# build the little-endian BGR->packed-int COLORREF value by hand.
function Get-RGBInt($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = Get-RGBInt 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = Get-RGBInt 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = Get-RGBInt 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = Get-RGBInt 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = Get-RGBInt 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = Get-RGBInt 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = Get-RGBInt 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = Get-RGBInt 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = Get-RGBInt 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = Get-RGBInt 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = Get-RGBInt 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = Get-RGBInt 0x95 0x4F 0x72   # folHlink
